$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Data for the 7 new people (rows 2-8 on Sheet1)
$emails = @(
    "sourabh.awasthi@capgemini.com",
    "sandipan.deb@capgemini.com",
    "biswaji.deb@capgemini.com",
    "debanjan.das@capgemini.com",
    "dhiraj.kajari@capgemini.com",
    "manoj-kumar.b.s@capgemini.com",
    "mayur.bhorkar@capgemini.com"
)

$passwords = @(
    '2dzo3m()F+Dl',
    'cg]6F)cxa}tE',
    'W@l3H%(O*Twr',
    '4VZBJZv+naBU',
    '3}{N?*W@D1u]',
    '&1IlV6S3rGH%',
    'p=nLFilY0!uh'
)

# --- Phase 1: literal values only, in the same order the data was originally typed in ---
# Rows 2 and 3 entered first (Sourabh Awasthi, Sandipan Deb)
$ws1.Range("D2").Value = $emails[0]
$ws1.Range("J2").Value = $passwords[0]
$ws1.Range("D3").Value = $emails[1]
$ws1.Range("J3").Value = $passwords[1]

# Then the lookup rows on Sheet2 for Sandipan Deb were filled in
$ws2.Range("C2").Value = "sandipan.deb"
$ws2.Range("I2").Value = "Capgemini"
$ws2.Range("C3").Value = "sandipan.deb"
$ws2.Range("I3").Value = "Capgemini"

# Then the remaining rows 4-8
for ($i = 2; $i -lt 7; $i++) {
    $r = $i + 2
    $ws1.Range("D$r").Value = $emails[$i]
    $ws1.Range("J$r").Value = $passwords[$i]
}

# --- Phase 2: formulas and remaining numeric/boolean values for every new row ---
for ($i = 0; $i -lt 7; $i++) {
    $r = $i + 2

    $ws1.Range("A$r").Formula = "=PROPER(IFERROR(LEFT(C$r,FIND(CHAR(46),C$r)-1),C$r))"
    $ws1.Range("B$r").Formula = '=IFERROR(PROPER(RIGHT(C' + $r + ',LEN(C' + $r + ')-FIND("@",SUBSTITUTE(C' + $r + ',".","@",((LEN(C' + $r + ')-LEN(SUBSTITUTE(C' + $r + ',".","")))/LEN("\")))))), "Unknown")'
    $ws1.Range("C$r").Formula = "=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D$r,FIND(CHAR(64),D$r)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))"
    $ws1.Range("E$r").Formula = "=LEFT(H$r,FIND(CHAR(46),H$r)-1)"
    $ws1.Range("F$r").Formula = "=CONCATENATE(" + '"ITPartner\"' + ",I$r)"
    $ws1.Range("H$r").Formula = "=RIGHT(D$r,LEN(D$r)-FIND(CHAR(64),D$r))"
    $ws1.Range("I$r").Formula = "=PROPER(E$r)"

    $ws1.Range("K$r").Value = 80
    $ws1.Range("M$r").Value = $true

    $ws1.Range("P$r").Formula = "=COUNTIF(D:D,D$r)"
}

Write-Host "applied"
